# Complete "Create a PivotTable and analyze your data"
# Adds two new "Title and Content" slides at the end of the deck:
#   slide24 - Excel动手实验室 – 数据透视表（创建）
#   slide25 - Excel动手实验室 – 数据分析表（分析）

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 24: "数据透视表（创建）" (Create a PivotTable)
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Add($p.Slides.Count + 1, 2)

$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Excel动手实验室 – 数据透视表（创建）"

$body1 = $s1.Shapes.Item(2).TextFrame.TextRange
$body1.Text = "创建数据透视表`r手动创建数据透视表`r字段`r行`r列`r值`r分组`r设计`r"

$body1.Paragraphs(3,1).IndentLevel = 2
$body1.Paragraphs(4,1).IndentLevel = 3
$body1.Paragraphs(5,1).IndentLevel = 3
$body1.Paragraphs(6,1).IndentLevel = 3
$body1.Paragraphs(7,1).IndentLevel = 2
$body1.Paragraphs(8,1).IndentLevel = 2

# ---------------------------------------------------------------------------
# Slide 25: "数据分析表（分析）" (Analyze your data)
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Add($p.Slides.Count + 1, 2)

$s2.Shapes.Item(1).TextFrame.TextRange.Text = "Excel动手实验室 – 数据分析表（分析）"

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "基本分析`r排序`r筛选`r汇总`r计算`r高级分析`r切片器`r日程表`r数据透视图"

$body2.Paragraphs(2,1).IndentLevel = 2
$body2.Paragraphs(3,1).IndentLevel = 2
$body2.Paragraphs(4,1).IndentLevel = 2
$body2.Paragraphs(5,1).IndentLevel = 2
$body2.Paragraphs(7,1).IndentLevel = 2
$body2.Paragraphs(8,1).IndentLevel = 2
$body2.Paragraphs(9,1).IndentLevel = 2
